$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are written as text to match the source data (prices/percentages
# are stored as formatted strings, not numbers). For values that look like plain
# numbers, force a Text number format before the write so Excel does not coerce
# them to a numeric type, then clear the format again so no stray style sticks.
$ws.Range("D2").Value = "67.826.13"
$ws.Range("E2").Value = "  +3.37%  "
$ws.Range("D3").Value = "3.329.87"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.78"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +6.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.46"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.65%  "
$ws.Range("D9").Value = "3.327.01"
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.179"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.581"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.45"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000276"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "642.63"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +12.71%  "
$ws.Range("D15").Value = "3.858.86"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.46"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.91%  "
$ws.Range("D17").Value = "67.971.13"
$ws.Range("E17").Value = "  +3.69%  "
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").Value = "3.331.53"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.72"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.93"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.900"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.69"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.71"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.61"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.10"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +9.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.58"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.66"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "594.43"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.95%  "
$ws.Range("D33").Value = "3.939.77"
$ws.Range("E33").Value = "  +5.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.95"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.56"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.104"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.74%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.82"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.27"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.02%  "
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.68"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "32.67"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").Value = "0.0₃0685"
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.338"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0416"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.52%  "
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.56"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.70%  "
$ws.Range("E50").Value = "  +10.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.16"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.37%  "
